$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.22%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.12%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.692"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.10%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08392"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.93%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.813"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.84%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'2.010"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.54%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.473"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.32%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.900"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.31%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9241"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.09%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1278"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.49%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1982"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.86%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09480"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.98%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03849"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.52%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1061"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.90%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001304"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.09%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-3.20%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.424"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.88%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'8.766"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.41%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1362"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-4.04%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-5.52%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04416"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.45%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001272"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.80%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004397"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.27%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001221"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-1.79%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-0.05%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02851"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.06%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05529"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.98%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007959"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.34%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1433"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.06%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-9.50%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-2.99%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01173"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.49%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006947"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.61%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003467"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'14.99%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.24%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").Style = "Normal"
